$d = $word.ActiveDocument

function Toggle-BoldRun($start, $end) {
    # Toggling Bold off then back on for a sub-range forces the engine to
    # materialize a distinct run boundary at (start,end) without altering
    # the final (bold) formatting value.
    $r = $d.Range($start, $end)
    $r.Bold = 0
    $r.Bold = 1
}

# The first paragraph currently reads:
#   "Supplementary Table 4 | Full-length insulator sequences."
# and must become:
#   "Supplementary Table 4. Full-length insulator sequences."
# with "Full-length insulator sequences." no longer bold while
# "Supplementary Table 4. " stays bold (and "Supplementary Table " gets
# split into three runs: "Supplementa" / "ry" / " Table ").

# Step 1: collapse " | " into ". " (text-only change; formatting for this
# span is untouched so the engine naturally merges the touched run with its
# neighbours -- that's fine, we re-split everything explicitly afterwards).
$sep = $d.Range(21, 24)
$sep.Text = ". "

# Step 2: re-establish explicit run boundaries now that the text is final.
# a) split "Supplementary Table " -> "Supplementa" | "ry" | " Table "
Toggle-BoldRun 11 13
# b) split "4. " -> "4" | "." | " "
Toggle-BoldRun 20 21
Toggle-BoldRun 21 22
Toggle-BoldRun 22 23

# Step 3: make "Full-length insulator sequences." non-bold (real formatting
# change, not just a toggle) -- this both creates the final run boundary and
# matches the target appearance.
$tail = $d.Range(23, 55)
$tail.Bold = 0
